$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.755.93"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.087.73"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.635"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.07"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.23"
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("D13").Value = "2.395.79"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "2.091.05"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "37.739.16"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.00"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.52"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.68"
$ws.Range("E26").Value = "  +8.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.82"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.57"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0636"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.62"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0236"
$ws.Range("E40").Value = "  +10.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.45"
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.74"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").Value = "1.450.30"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.06"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.22"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "2.279.55"
$ws.Range("E51").Value = "  +0.91%  "
